$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Range("A23").Value = 41193
$ws.Range("B23").Value = 2.25
$ws.Range("D23").Value = "Manual continued"

# Row 24
$ws.Range("A24").Value = 41194
$ws.Range("B24").Value = 5
$ws.Range("D24").Value = "Manual continued"

# Row 25
$ws.Range("A25").Value = 41197
$ws.Range("B25").Value = 1
$ws.Range("D25").Value = "Fix: Bad specification of ALL events - now timer events are still an OR condition"

# Update selection to match the diff
$ws.Range("E25").Select()
